$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# Each entry: (row index [1-based], cell index [1-based], new text)
# Row 10 = Panama, Row 15 = Peru, Row 16 = Brazil, Row 17 = Chile,
# Row 18 = Bolivia, Row 19 = Colombia
# Cell 3 = Confirmed Cases, Cell 4 = Cummulative Incidence (per million)

$updates = @(
    @(10, 3, "7"),
    @(10, 4, "1.59"),
    @(15, 3, "1,068"),
    @(15, 4, "31.37"),
    @(16, 3, "3,756"),
    @(16, 4, "17.44"),
    @(17, 3, "207"),
    @(17, 4, "10.56"),
    @(18, 3, "42"),
    @(18, 4, "3.44"),
    @(19, 3, "164"),
    @(19, 4, "3.16")
)

foreach ($u in $updates) {
    $rowIdx = $u[0]
    $colIdx = $u[1]
    $newText = $u[2]
    $cell = $tbl.Rows($rowIdx).Cells($colIdx)
    $cell.Range.Text = $newText
}
